$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(305487936, Avihai  Kipnis: 4,4)"
$ws.Range("B1").Value = "(313227928, Aviv  Levi: -7,-8)"
$ws.Range("C1").Value = "(205807308, Sariel  Basis: 7,-4)"
$ws.Range("D1").Value = "(315891549, Raz  Halaby: 1,7)"
$ws.Range("E1").Value = "(315060103, Dan  Mshelh: 0,-3)"
$ws.Range("F1").Value = "(313925141, Elad   Amer: -8,2)"
$ws.Range("G1").Value = "(326598423, Ron Cohen: -1,-6)"

$ws.Range("A3").Value = "cost: 419.9321605734631"
$ws.Range("A4").Value = "time: 56.4188800819233"
